$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255, shifting existing rows 255..308 down to 256..309
$ws.Rows.Item(255).Insert()

# Populate the new row 255 with the new record's data
$ws.Cells.Item(255, 1).Value = 4
$ws.Cells.Item(255, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(255, 3).Value = "Los Lagos"
$ws.Cells.Item(255, 4).Value = 44711
$ws.Cells.Item(255, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(255, 5).Value = 10
$ws.Cells.Item(255, 6).Value = 100112045
$ws.Cells.Item(255, 7).Value = "Zapallo"
$ws.Cells.Item(255, 8).Value = "Paine"
$ws.Cells.Item(255, 9).Value = "1a (cosecha)"
$ws.Cells.Item(255, 10).Value = 500
$ws.Cells.Item(255, 11).Value = 500
$ws.Cells.Item(255, 12).Value = 500
$ws.Cells.Item(255, 13).Value = 500
$ws.Cells.Item(255, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(255, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(255, 16).Value = 500
$ws.Cells.Item(255, 17).Value = 1
$ws.Cells.Item(255, 18).Value = "Hortaliza"
